# Applies the "Group contract" edits described in the commit diff.
#
# Strategy: for every paragraph whose run layout changes (either because a
# name needs spell-check run-splitting via <w:proofErr/>, or because a
# trailing sentence is replaced/extended with several new runs), we rebuild
# the *entire* paragraph's run content and push it back with
# Range.InsertXML. InsertXML replaces exactly the contents of the range it
# is called on, so by targeting the paragraph's full range (excluding the
# trailing paragraph mark) and wrapping the replacement runs in a <w:p>
# element, the paragraph's own <w:pPr>/identity is preserved while its run
# children are swapped out atomically.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-ParaXml {
    param(
        [int]$ParaIndex,
        [string]$InnerXml
    )
    $para = $d.Paragraphs.Item($ParaIndex)
    $full = $para.Range
    $target = $d.Range($full.Start, $full.End - 1)
    $xml = "<w:p $wNs>$InnerXml</w:p>"
    $target.InsertXML($xml)
}

# 1) "Tor Oveland Gikling" -> split with spell-check markers around the
#    two surnames.
Set-ParaXml 4 '<w:r><w:t xml:space="preserve">Tor </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Oveland</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Gikling</w:t></w:r><w:proofErr w:type="spellEnd"/>'

# 2) "Anders M. H. Frostrud" -> split off the surname with spell-check markers.
Set-ParaXml 5 '<w:r><w:t xml:space="preserve">Anders M. H. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Frostrud</w:t></w:r><w:proofErr w:type="spellEnd"/>'

# 3) "Thomas Ystenes" -> split off the surname with spell-check markers.
Set-ParaXml 6 '<w:r><w:t xml:space="preserve">Thomas </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Ystenes</w:t></w:r><w:proofErr w:type="spellEnd"/>'

# 4) Team leader bullet: drop trailing space after "represents the group"
#    and append who fills that role.
Set-ParaXml 10 '<w:r><w:t>T</w:t></w:r><w:r w:rsidR="00AB11A5"><w:t>eam leader</w:t></w:r><w:r w:rsidR="00E94937"><w:t>, c</w:t></w:r><w:r w:rsidR="00AB11A5"><w:t xml:space="preserve">ommunicator </w:t></w:r><w:r w:rsidR="00E94937"><w:t xml:space="preserve">- </w:t></w:r><w:r w:rsidR="00AB11A5"><w:t>calls in for meetings</w:t></w:r><w:r w:rsidR="00E94937"><w:t xml:space="preserve">, brings the team together, </w:t></w:r><w:r w:rsidR="00860998"><w:t>represents the group</w:t></w:r><w:r><w:t xml:space="preserve">. (Tor </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Oveland</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Gikling</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>).</w:t></w:r>'

# 5) Documentation manager bullet: append who fills that role.
Set-ParaXml 11 '<w:r><w:t>D</w:t></w:r><w:r w:rsidR="00C11409"><w:t>ocumentation</w:t></w:r><w:r w:rsidR="00AB11A5"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00C11409"><w:t xml:space="preserve">manager </w:t></w:r><w:r w:rsidR="00AB11A5"><w:t>(</w:t></w:r><w:r w:rsidR="00C11409"><w:t>makes sure all the documents are in place, does not write all him/her-self, but coordinates and makes sure that the work is done in a team)</w:t></w:r><w:r><w:t xml:space="preserve"> (Thomas </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Ystenes</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>).</w:t></w:r>'

# 6) Quality control bullet: add a trailing period.
Set-ParaXml 12 '<w:r><w:t>Quality control</w:t></w:r><w:r w:rsidR="0087562A"><w:t xml:space="preserve"> – checks code, contributes with testing, </w:t></w:r><w:r w:rsidR="00066E0D"><w:t>reads through documents, checks for errors</w:t></w:r><w:r><w:t>.</w:t></w:r>'

# 7) Architect bullet: add a trailing period.
Set-ParaXml 13 '<w:r><w:t>Architect – makes sure that code architecture is clean</w:t></w:r><w:r><w:t>.</w:t></w:r>'

# 8) Meetings bullet: replace the description with the new meeting schedule.
Set-ParaXml 16 '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Meetings</w:t></w:r><w:r w:rsidR="003E68A7"><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t xml:space="preserve">8:15-11 </w:t></w:r><w:r><w:t>Wednesday</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>12-14 Friday.</w:t></w:r>'

# 9) Checkpoints bullet: replace the description with the new checkpoint plan.
Set-ParaXml 18 '<w:r w:rsidRPr="00800FCB"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Checkpoints</w:t></w:r><w:r w:rsidR="005849AD"><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t>30 minutes for weekly checkpoints on Friday</w:t></w:r><w:r><w:t>.</w:t></w:r>'

# 10) Absence notification bullet: tweak wording and add advance-notice rule.
Set-ParaXml 19 '<w:r w:rsidRPr="00800FCB"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Absence notification</w:t></w:r><w:r w:rsidR="00800FCB"><w:t>. I</w:t></w:r><w:r><w:t>f</w:t></w:r><w:r><w:t xml:space="preserve"> a group member is late or </w:t></w:r><w:r><w:t>cannot</w:t></w:r><w:r><w:t xml:space="preserve"> attend, a message to notify the other group members shall be sent</w:t></w:r><w:r><w:t xml:space="preserve"> at least 2 hours prior</w:t></w:r><w:r><w:t>.</w:t></w:r>'

# 11) Documentation (procedures) bullet: "files, and" -> "files and".
Set-ParaXml 20 '<w:r w:rsidRPr="00800FCB"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Documentation</w:t></w:r><w:r w:rsidR="00800FCB"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t xml:space="preserve">We use a GitHub repository to manage documents and </w:t></w:r><w:r><w:t>files and</w:t></w:r><w:r><w:t xml:space="preserve"> keep track on changes and updates.</w:t></w:r><w:r w:rsidR="00235D00"><w:t xml:space="preserve"> </w:t></w:r>'
